$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (id) was reformatted as Text so large numeric ids (e.g. employee
# numbers) display / sort exactly as typed instead of being treated as numbers.
$ws.Range("B2:B15").NumberFormat = "@"
$ws.Range("B16:B17").NumberFormat = "@"
$ws.Range("B16:B17").WrapText = $true

# --- D17 previously stored a literal boolean; make it consistent with every
# other row in column D, which stores the result of a FALSE() formula.
$ws.Range("D17").Formula = "=FALSE()"

# --- Add the new user row for the additional water-connection role.
$ws.Range("A18").Value = "seniorAssistant1"
$ws.Range("B18").Value = 944176
$ws.Range("B18").NumberFormat = "@"
$ws.Range("C18").Value = "kurnool_eGov@123"
$ws.Range("D18").Formula = "=FALSE()"

$ws.Hyperlinks.Add($ws.Range("C18"), "mailto:kurnool_eGov@123", "", "", "kurnool_eGov@123") | Out-Null
